# Update the "Metadata" sheet (FHIR StructureDefinition metadata table).
$wb = $excel.ActiveWorkbook
$metadata = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$metadata.Range("B3").Value = "6.0.0"

# Date: refreshed publication timestamp
$metadata.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher used to be blank ("Contact" row carried no real value) -- now populated.
$metadata.Range("B9").Value = "Alvearie Team"

# The old row 10 was a "Contact" property with placeholder text; it is replaced
# with a "Jurisdiction" property describing the applicable jurisdiction.
$metadata.Range("A10").Value = "Jurisdiction"
$metadata.Range("B10").Value = "United States of America"

# The old row 11 duplicated row 10 ("Contact" / "No display for ContactDetail")
# and is removed entirely, shifting all subsequent rows up by one.
$metadata.Rows.Item(11).Delete()

# Update the "Elements" sheet: the root Extension row's Short/Definition
# columns now mirror the StructureDefinition's Title/Description instead of
# the generic placeholder text.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Mental Health Ambulatory Coverage Indicator"
$elements.Range("L2").Value = "Indicates whether the member has mental health ambulatory benefit coverage: Y or N. This finer granularity of MHSA benefit coverage may be used in HEDIS reporting."

# Column K (Short) grows to fit the longer text (best-fit width).
$elements.Columns.Item(11).ColumnWidth = 42.3
